$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.492.44'
$ws.Range('E2').Value = '  -2.59%  '
$ws.Range('D3').Value = '1.806.08'
$ws.Range('E3').Value = '  -2.45%  '
$ws.Range('E4').Value = '  +0.75%  '
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4562'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3663'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07127'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8787'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07739'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.37'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.96%  '
$ws.Range('D13').Value = '1.811.79'
$ws.Range('E13').Value = '  -7.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.274'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.361'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '86.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008577'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.69%  '
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').Value = '26.532.62'
$ws.Range('E20').Value = '  -2.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.981'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('E23').Value = '  -0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.985'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.90'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.044'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '112.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.844'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08664'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.041'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7300'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.447'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('E34').Value = '  -4.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.006'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.541'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.080'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01934'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05106'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.894'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.947'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5012'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('E43').Value = '  -3.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.136'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.90%  '
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4604'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.72%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.942'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.13'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.592'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05990'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.99'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.82%  '
